$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new email address (a "gvu." sub-domain was inserted right
# after the "@" of each giáo vụ's e-mail address).
$emails = [ordered]@{
  "C2" = "nguyenvanb@gvu.soict.hust.edu.vn"
  "C3" = "lythic@gvu.spkt.hust.edu.vn"
  "C4" = "lethidc@gvu.nn.hust.edu.vn"
  "C5" = "tranvane@gvu.dtvt.hust.edu.vn"
  "C6" = "daothif@gvu.dktdh.hust.edu.vn"
}

foreach ($addr in $emails.Keys) {
  $email = $emails[$addr]
  # Update the visible text first ...
  $ws.Range($addr).Value = $email
  # ... then turn the cell into a live mailto: hyperlink (this also applies
  # the built-in "Hyperlink" cell style: underline + theme color).
  $ws.Hyperlinks.Add($ws.Range($addr), "mailto:$email") | Out-Null
}

# Selection moved from E5 to C7 in the saved view state.
$ws.Range("C7").Select() | Out-Null
